$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "246.27"
Set-TextValue "D4" "5.319"
Set-TextValue "D5" "0.05875"
Set-TextValue "D6" "3.394"
Set-TextValue "D8" "0.8127"
Set-TextValue "D9" "0.9560"
Set-TextValue "D10" "0.1413"
Set-TextValue "D11" "0.03662"
Set-TextValue "D12" "0.07327"
Set-TextValue "D13" "0.03054"
Set-TextValue "D15" "0.09402"
Set-TextValue "D16" "0.001601"
Set-TextValue "D18" "0.0005905"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006170"
Set-TextValue "D21" "0.0009864"
Set-TextValue "D22" "0.00009707"
Set-TextValue "D23" "3.686"
Set-TextValue "D26" "0.1285"
Set-TextValue "D27" "0.0002474"
Set-TextValue "D40" "0.03896"
Set-TextValue "D41" "0.006755"
Set-TextValue "D43" "0.003002"
Set-TextValue "D44" "0.005917"
Set-TextValue "D45" "0.00005672"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.6521"
Set-TextValue "D48" "0.06606"
Set-TextValue "E48" "47BOLOBOLO"
Set-TextValue "D50" "0.01011"

Write-Output "Done"
